$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, copying the formatting used by the
# other header cells (e.g. G1 - bold, centered, bordered style).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add corresponding data value in H2
$ws.Range("H2").Value = 0
